$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the "addWithCarryS" documentation row. Its function-signature cell
# (C2) was missing the "ShiftDirection" parameter that the implementation
# actually takes, and the sheet had no checkmark in the (new) ShiftDirection
# column (J) for this function. Bring the row up to date:

# 1) Update the function signature to include "-> ShiftDirection"
$ws.Range("C2").Value = "RegisterID -> RegisterID -> Operand -> MachineState -> bool -> bool -> ShiftDirection -> MachineState"

# 2) Mark that addWithCarryS does take a ShiftDirection argument (column J)
$ws.Range("J2").Value = "$([char]0x2713)"
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").WrapText = $true

# 3) Leave the selection on the cell that was edited
[void]$ws.Range("C2").Select()
